# Runtime update: refresh referee penalty-minute stats and the as_of_utc timestamp
# for every row on the "Главные" and "Линейные" sheets, per the 2025-10-29 data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Главные")

$ws.Range("C2").Value = 18
$ws.Range("D2").Value = 452
$ws.Range("E2").Value = 194
$ws.Range("F2").Value = 258
$ws.Range("G2").Value = 25.11
$ws.Range("H2").Value = 10.78
$ws.Range("I2").Value = 14.33
$ws.Range("J2").Value = 82
$ws.Range("K2").Value = 99
$ws.Range("O2").Value = 1
$ws.Range("Q2").Value = 2
$ws.Range("AA2").Value = "2025-10-29 07:08:09"
$ws.Range("C3").Value = 18
$ws.Range("D3").Value = 325
$ws.Range("E3").Value = 147
$ws.Range("F3").Value = 178
$ws.Range("G3").Value = 18.06
$ws.Range("H3").Value = 8.17
$ws.Range("I3").Value = 9.890000000000001
$ws.Range("J3").Value = 71
$ws.Range("K3").Value = 74
$ws.Range("AA3").Value = "2025-10-29 07:08:09"
$ws.Range("C4").Value = 13
$ws.Range("D4").Value = 228
$ws.Range("E4").Value = 97
$ws.Range("F4").Value = 131
$ws.Range("G4").Value = 17.54
$ws.Range("H4").Value = 7.46
$ws.Range("J4").Value = 46
$ws.Range("K4").Value = 53
$ws.Range("AA4").Value = "2025-10-29 07:08:09"
$ws.Range("AA5").Value = "2025-10-29 07:08:09"
$ws.Range("C6").Value = 18
$ws.Range("D6").Value = 351
$ws.Range("E6").Value = 153
$ws.Range("F6").Value = 198
$ws.Range("G6").Value = 19.5
$ws.Range("H6").Value = 8.5
$ws.Range("I6").Value = 11
$ws.Range("J6").Value = 69
$ws.Range("K6").Value = 84
$ws.Range("AA6").Value = "2025-10-29 07:08:09"
$ws.Range("AA7").Value = "2025-10-29 07:08:09"
$ws.Range("C8").Value = 17
$ws.Range("D8").Value = 324
$ws.Range("E8").Value = 151
$ws.Range("F8").Value = 173
$ws.Range("G8").Value = 19.06
$ws.Range("H8").Value = 8.880000000000001
$ws.Range("I8").Value = 10.18
$ws.Range("J8").Value = 68
$ws.Range("K8").Value = 79
$ws.Range("AA8").Value = "2025-10-29 07:08:09"
$ws.Range("AA9").Value = "2025-10-29 07:08:09"
$ws.Range("AA10").Value = "2025-10-29 07:08:09"
$ws.Range("AA11").Value = "2025-10-29 07:08:09"
$ws.Range("AA12").Value = "2025-10-29 07:08:09"
$ws.Range("AA13").Value = "2025-10-29 07:08:09"
$ws.Range("C14").Value = 12
$ws.Range("D14").Value = 149
$ws.Range("E14").Value = 72
$ws.Range("F14").Value = 77
$ws.Range("G14").Value = 12.42
$ws.Range("H14").Value = 6
$ws.Range("I14").Value = 6.42
$ws.Range("J14").Value = 36
$ws.Range("K14").Value = 36
$ws.Range("AA14").Value = "2025-10-29 07:08:09"
$ws.Range("AA15").Value = "2025-10-29 07:08:09"
$ws.Range("C16").Value = 19
$ws.Range("D16").Value = 375
$ws.Range("E16").Value = 178
$ws.Range("F16").Value = 197
$ws.Range("G16").Value = 19.74
$ws.Range("H16").Value = 9.369999999999999
$ws.Range("I16").Value = 10.37
$ws.Range("J16").Value = 69
$ws.Range("K16").Value = 71
$ws.Range("AA16").Value = "2025-10-29 07:08:09"
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 192
$ws.Range("E17").Value = 64
$ws.Range("F17").Value = 128
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 5.33
$ws.Range("I17").Value = 10.67
$ws.Range("J17").Value = 32
$ws.Range("K17").Value = 49
$ws.Range("O17").Value = 1
$ws.Range("Q17").Value = 1
$ws.Range("AA17").Value = "2025-10-29 07:08:09"
$ws.Range("AA18").Value = "2025-10-29 07:08:09"
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 272
$ws.Range("E19").Value = 122
$ws.Range("F19").Value = 150
$ws.Range("G19").Value = 18.13
$ws.Range("H19").Value = 8.130000000000001
$ws.Range("I19").Value = 10
$ws.Range("J19").Value = 56
$ws.Range("M19").Value = 4
$ws.Range("V19").Value = 8
$ws.Range("AA19").Value = "2025-10-29 07:08:09"
$ws.Range("AA20").Value = "2025-10-29 07:08:09"
$ws.Range("AA21").Value = "2025-10-29 07:08:09"
$ws.Range("AA22").Value = "2025-10-29 07:08:09"
$ws.Range("C23").Value = 11
$ws.Range("D23").Value = 138
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 88
$ws.Range("G23").Value = 12.55
$ws.Range("H23").Value = 4.55
$ws.Range("I23").Value = 8
$ws.Range("J23").Value = 25
$ws.Range("K23").Value = 34
$ws.Range("AA23").Value = "2025-10-29 07:08:09"
$ws.Range("AA24").Value = "2025-10-29 07:08:09"
$ws.Range("AA25").Value = "2025-10-29 07:08:09"
$ws.Range("AA26").Value = "2025-10-29 07:08:09"

$ws = $wb.Worksheets.Item("Линейные")

$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 208
$ws.Range("E2").Value = 74
$ws.Range("F2").Value = 134
$ws.Range("G2").Value = 18.91
$ws.Range("H2").Value = 6.73
$ws.Range("I2").Value = 12.18
$ws.Range("J2").Value = 37
$ws.Range("K2").Value = 47
$ws.Range("AA2").Value = "2025-10-29 07:08:09"
$ws.Range("C3").Value = 17
$ws.Range("D3").Value = 248
$ws.Range("E3").Value = 118
$ws.Range("F3").Value = 130
$ws.Range("G3").Value = 14.59
$ws.Range("H3").Value = 6.94
$ws.Range("I3").Value = 7.65
$ws.Range("J3").Value = 59
$ws.Range("M3").Value = 2
$ws.Range("V3").Value = 10
$ws.Range("AA3").Value = "2025-10-29 07:08:09"
$ws.Range("AA4").Value = "2025-10-29 07:08:09"
$ws.Range("AA5").Value = "2025-10-29 07:08:09"
$ws.Range("AA6").Value = "2025-10-29 07:08:09"
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 175
$ws.Range("E7").Value = 62
$ws.Range("F7").Value = 113
$ws.Range("G7").Value = 17.5
$ws.Range("H7").Value = 6.2
$ws.Range("I7").Value = 11.3
$ws.Range("J7").Value = 31
$ws.Range("K7").Value = 34
$ws.Range("AA7").Value = "2025-10-29 07:08:09"
$ws.Range("AA8").Value = "2025-10-29 07:08:09"
$ws.Range("AA9").Value = "2025-10-29 07:08:09"
$ws.Range("AA10").Value = "2025-10-29 07:08:09"
$ws.Range("C11").Value = 16
$ws.Range("D11").Value = 298
$ws.Range("E11").Value = 143
$ws.Range("F11").Value = 155
$ws.Range("G11").Value = 18.63
$ws.Range("H11").Value = 8.94
$ws.Range("I11").Value = 9.69
$ws.Range("J11").Value = 64
$ws.Range("K11").Value = 70
$ws.Range("AA11").Value = "2025-10-29 07:08:09"
$ws.Range("AA12").Value = "2025-10-29 07:08:09"
$ws.Range("C13").Value = 16
$ws.Range("D13").Value = 347
$ws.Range("E13").Value = 175
$ws.Range("F13").Value = 172
$ws.Range("G13").Value = 21.69
$ws.Range("H13").Value = 10.94
$ws.Range("I13").Value = 10.75
$ws.Range("J13").Value = 65
$ws.Range("K13").Value = 66
$ws.Range("AA13").Value = "2025-10-29 07:08:09"
$ws.Range("AA14").Value = "2025-10-29 07:08:09"
$ws.Range("AA15").Value = "2025-10-29 07:08:09"
$ws.Range("C16").Value = 19
$ws.Range("D16").Value = 359
$ws.Range("E16").Value = 170
$ws.Range("F16").Value = 189
$ws.Range("G16").Value = 18.89
$ws.Range("H16").Value = 8.949999999999999
$ws.Range("I16").Value = 9.949999999999999
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = 77
$ws.Range("AA16").Value = "2025-10-29 07:08:09"
$ws.Range("AA17").Value = "2025-10-29 07:08:09"
$ws.Range("C18").Value = 14
$ws.Range("D18").Value = 230
$ws.Range("E18").Value = 117
$ws.Range("F18").Value = 113
$ws.Range("G18").Value = 16.43
$ws.Range("H18").Value = 8.359999999999999
$ws.Range("I18").Value = 8.07
$ws.Range("J18").Value = 56
$ws.Range("K18").Value = 54
$ws.Range("AA18").Value = "2025-10-29 07:08:09"
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 450
$ws.Range("E19").Value = 197
$ws.Range("F19").Value = 253
$ws.Range("G19").Value = 21.43
$ws.Range("H19").Value = 9.380000000000001
$ws.Range("J19").Value = 91
$ws.Range("K19").Value = 104
$ws.Range("AA19").Value = "2025-10-29 07:08:09"
$ws.Range("AA20").Value = "2025-10-29 07:08:09"
$ws.Range("AA21").Value = "2025-10-29 07:08:09"
$ws.Range("C22").Value = 20
$ws.Range("D22").Value = 390
$ws.Range("E22").Value = 155
$ws.Range("F22").Value = 235
$ws.Range("G22").Value = 19.5
$ws.Range("H22").Value = 7.75
$ws.Range("I22").Value = 11.75
$ws.Range("J22").Value = 70
$ws.Range("K22").Value = 90
$ws.Range("O22").Value = 1
$ws.Range("Q22").Value = 2
$ws.Range("AA22").Value = "2025-10-29 07:08:09"
$ws.Range("AA23").Value = "2025-10-29 07:08:09"
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 382
$ws.Range("E24").Value = 157
$ws.Range("F24").Value = 225
$ws.Range("G24").Value = 23.88
$ws.Range("H24").Value = 9.81
$ws.Range("I24").Value = 14.06
$ws.Range("J24").Value = 56
$ws.Range("K24").Value = 60
$ws.Range("O24").Value = 1
$ws.Range("Q24").Value = 3
$ws.Range("AA24").Value = "2025-10-29 07:08:09"
$ws.Range("AA25").Value = "2025-10-29 07:08:09"
$ws.Range("AA26").Value = "2025-10-29 07:08:09"
